$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2121212121212121
$ws.Range("C2").Value = 0.4915824915824916
$ws.Range("J2").Value = 0.0202020202020202
$ws.Range("P2").Value = 0.1717171717171717
$ws.Range("S2").Value = 0.1043771043771044
$ws.Range("B3").Value = 0.0131578947368421
$ws.Range("C3").Value = 0.02631578947368421
$ws.Range("J3").Value = 0.04605263157894737
$ws.Range("P3").Value = 0.7105263157894737
$ws.Range("S3").Value = 0.2039473684210526
$ws.Range("J4").Value = 0.06451612903225806
$ws.Range("P4").Value = 0.5806451612903226
$ws.Range("S4").Value = 0.3548387096774194
$ws.Range("B6").Value = 0.07425742574257425
$ws.Range("F6").Value = 0.07425742574257425
$ws.Range("J6").Value = 0.2425742574257426
$ws.Range("O6").Value = 0.02475247524752475
$ws.Range("Q6").Value = 0.1930693069306931
$ws.Range("R6").Value = 0.06435643564356436
$ws.Range("S6").Value = 0.3267326732673267
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.01169590643274854
$ws.Range("F7").Value = 0.04678362573099415
$ws.Range("J7").Value = 0.1228070175438596
$ws.Range("O7").Value = 0.03508771929824561
$ws.Range("Q7").Value = 0.1871345029239766
$ws.Range("R7").Value = 0.0935672514619883
$ws.Range("S7").Value = 0.391812865497076
$ws.Range("B8").Value = 0.1007751937984496
$ws.Range("D8").Value = 0.01808785529715762
$ws.Range("F8").Value = 0.06718346253229975
$ws.Range("J8").Value = 0.1291989664082687
$ws.Range("O8").Value = 0.02325581395348837
$ws.Range("Q8").Value = 0.165374677002584
$ws.Range("R8").Value = 0.06976744186046512
$ws.Range("S8").Value = 0.4263565891472868
$ws.Range("B9").Value = 0.1259259259259259
$ws.Range("D9").Value = 0.01481481481481482
$ws.Range("F9").Value = 0.08148148148148149
$ws.Range("J9").Value = 0.1481481481481481
$ws.Range("O9").Value = 0.007407407407407408
$ws.Range("Q9").Value = 0.1259259259259259
$ws.Range("R9").Value = 0.06666666666666667
$ws.Range("S9").Value = 0.4296296296296296
$ws.Range("B10").Value = 0.1071149335418295
$ws.Range("D10").Value = 0.01720093823299453
$ws.Range("E10").Value = 0.002345582486317436
$ws.Range("F10").Value = 0.06098514464425332
$ws.Range("J10").Value = 0.1532447224394058
$ws.Range("O10").Value = 0.02345582486317436
$ws.Range("Q10").Value = 0.2111024237685692
$ws.Range("R10").Value = 0.07193119624706802
$ws.Range("S10").Value = 0.3526192337763878
$ws.Range("G11").Value = 0.1503759398496241
$ws.Range("J11").Value = 0.09398496240601503
$ws.Range("K11").Value = 0.1804511278195489
$ws.Range("L11").Value = 0.5601503759398496
$ws.Range("S11").Value = 0.01503759398496241
$ws.Range("G12").Value = 0.7243589743589743
$ws.Range("J12").Value = 0.1794871794871795
$ws.Range("K12").Value = 0.01282051282051282
$ws.Range("L12").Value = 0.04487179487179487
$ws.Range("S12").Value = 0.03846153846153846
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2142857142857143
$ws.Range("S13").Value = 0.07142857142857142
$ws.Range("F15").Value = 0.02173913043478261
$ws.Range("H15").Value = 0.1358695652173913
$ws.Range("I15").Value = 0.03260869565217391
$ws.Range("J15").Value = 0.4239130434782609
$ws.Range("K15").Value = 0.07608695652173914
$ws.Range("M15").Value = 0.0108695652173913
$ws.Range("O15").Value = 0.02717391304347826
$ws.Range("S15").Value = 0.2717391304347826
$ws.Range("F16").Value = 0.05617977528089887
$ws.Range("H16").Value = 0.1348314606741573
$ws.Range("I16").Value = 0.06741573033707865
$ws.Range("J16").Value = 0.4606741573033708
$ws.Range("K16").Value = 0.1123595505617977
$ws.Range("M16").Value = 0.01123595505617977
$ws.Range("O16").Value = 0.02247191011235955
$ws.Range("S16").Value = 0.1348314606741573
$ws.Range("F17").Value = 0.01666666666666667
$ws.Range("H17").Value = 0.1666666666666667
$ws.Range("I17").Value = 0.08095238095238096
$ws.Range("J17").Value = 0.4523809523809524
$ws.Range("K17").Value = 0.09523809523809523
$ws.Range("M17").Value = 0.0119047619047619
$ws.Range("O17").Value = 0.05238095238095238
$ws.Range("S17").Value = 0.1238095238095238
$ws.Range("F18").Value = 0.03267973856209151
$ws.Range("H18").Value = 0.1568627450980392
$ws.Range("I18").Value = 0.07843137254901961
$ws.Range("J18").Value = 0.4705882352941176
$ws.Range("K18").Value = 0.1176470588235294
$ws.Range("M18").Value = 0.006535947712418301
$ws.Range("O18").Value = 0.05228758169934641
$ws.Range("S18").Value = 0.08496732026143791
$ws.Range("F19").Value = 0.008695652173913044
$ws.Range("H19").Value = 0.22
$ws.Range("I19").Value = 0.06086956521739131
$ws.Range("J19").Value = 0.4060869565217392
$ws.Range("K19").Value = 0.09826086956521739
$ws.Range("M19").Value = 0.01565217391304348
$ws.Range("N19").Value = 0.001739130434782609
$ws.Range("O19").Value = 0.06347826086956522
$ws.Range("S19").Value = 0.1252173913043478
